# The deck ships with two theme parts:
#   ppt/theme/theme1.xml  -> "Office Theme" colour scheme (used only by the
#                             notes/handout masters)
#   ppt/theme/theme2.xml  -> "Integral"     colour scheme (used by the one
#                             slide master, i.e. what is actually visible on
#                             the slides, and referenced directly from
#                             presentation.xml)
#
# The authored change swaps the *content* of those two theme parts: the
# slide master's theme becomes the stock "Office Theme" palette, while the
# (otherwise invisible) notes-master theme becomes the old "Integral"
# palette. Font scheme / format scheme are identical between the two parts,
# so only the twelve colour-scheme slots (dk1, lt1, dk2, lt2, accent1-6,
# hlink, folHlink) actually change value.
#
# Re-colour the presentation's (single, shared) theme through the standard
# PowerPoint colour-scheme object model so it ends up holding the "Office
# Theme" palette that theme2.xml should now contain.

$p = $ppt.ActivePresentation

$officeTheme = @(
    0,          # dk1      000000
    16777215,   # lt1      FFFFFF
    6968388,    # dk2      44546A
    15132391,   # lt2      E7E6E6
    13998939,   # accent1  5B9BD5
    3243501,    # accent2  ED7D31
    10855845,   # accent3  A5A5A5
    49407,      # accent4  FFC000
    12874308,   # accent5  4472C4
    4697456,    # accent6  70AD47
    12673797,   # hlink    0563C1
    7491477     # folHlink 954F72
)

$slide = $p.Slides.Item(1)
$scheme = $slide.ThemeColorScheme

for ($i = 1; $i -le $scheme.Count; $i++) {
    $scheme.Item($i).RGB = $officeTheme[$i - 1]
}

Write-Output "Re-applied Office Theme colour scheme to the presentation theme."
